# Working on replacing code chunks with embed shortcodes.
#
# Inserts a new "Body Text" styled paragraph immediately after each of the
# three chunk-output images (the "First Paragraph" style paragraphs that
# hold the inline drawings produced by the R/Python/Julia code chunks), and
# immediately before the SourceCode paragraph that follows it.
#
# NOTE: this runtime loses an inline drawing if a paragraph mark is created
# directly on/adjacent to the drawing-bearing paragraph's own Range and then
# text is written into it. To avoid corrupting the pictures, new paragraphs
# are built in a safe location (appended at the very end of the story, which
# never touches a drawing) and then moved into place with Cut/Paste, which
# is the pattern found to reliably preserve the surrounding content.

$d = $word.ActiveDocument

function Insert-BodyTextParagraphAfter($paragraphIndex, $text) {
    # Build the paragraph somewhere neutral (document end) so the drawing in
    # $paragraphIndex is never touched directly.
    $null = $d.Paragraphs.Add()
    $lastIndex = $d.Paragraphs.Count
    $scratch = $d.Paragraphs.Item($lastIndex)
    $scratch.Range.InsertAfter($text)

    # Cut the freshly typed paragraph (text + its paragraph mark) ...
    $scratch.Range.Cut()

    # ... and paste it right after the target paragraph.
    $target = $d.Paragraphs.Item($paragraphIndex)
    $pastePoint = $d.Range($target.Range.End, $target.Range.End)
    $pastePoint.Paste()

    # The pasted paragraph is now the one right after the target; style it.
    $newParagraph = $d.Paragraphs.Item($paragraphIndex + 1)
    $newParagraph.Style = "Body Text"
}

# Paragraph 18 = image for the R/ggpairs chunk (unnamed-chunk-1)
Insert-BodyTextParagraphAfter 18 "Another test"

# After the first insertion everything shifted down by one, so the image
# for the Python chunk (unnamed-chunk-2) is now paragraph 21.
Insert-BodyTextParagraphAfter 21 "Will this be deleted?"

# After the second insertion everything shifted down by one more, so the
# image for the Julia chunk (unnamed-chunk-3) is now paragraph 24.
Insert-BodyTextParagraphAfter 24 "Test"
